$wb = $excel.ActiveWorkbook

# --- Hoja1: update the "Conversión del día" note with the new rates ---
$ws1 = $wb.Worksheets.Item("Hoja1")

$text = @'
Conversión del día 💰
✅ Dólar paralelo: 68

Binance
✅ 1000 Bs = 14.85 = 61995.25 pesos
✅ 61995.25 pesos = 14.77 = 977.45 Bs

Promedio competencia
✅ Tasa pesos: 20
✅ Tasa Bs: 20
✅ % Ganancia: 20%
'@

$ws1.Range("A1").Value = $text

# --- tasas: refresh the N10/O10/N12/O12 rate figures ---
$ws2 = $wb.Worksheets.Item("tasas")

$ws2.Range("N10").Value = 67.36
$ws2.Range("O10").Value = 4176
$ws2.Range("N12").Value = 4195.99
$ws2.Range("O12").Value = 66.15600000000001
